# Falcon PoC - DataPipeline.pptx edits
# - Bump the three cached "datetimeFigureOut" field captions from 9/6/15 to 9/8/15
#   (Blank slide layout, Handout Master, Notes Master).
# - Merge the two "Data "/"Pipeline Demo" title runs on slide 1 into one run.
# - Remove the stray empty "Text Placeholder 2" (idx=10) shapes on slides 14 & 15.
# - Re-word the Subtitle on slide 14 from "Lets do Hadoop " to "We do Hadoop ".

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Date placeholder on the "Blank" slide layout (used by slide 2, etc.)
# ---------------------------------------------------------------------------
$blankLayout = $p.Slides.Item(2).CustomLayout
$blankLayout.Shapes.Item(1).TextFrame.TextRange.Text = "9/8/15"

# ---------------------------------------------------------------------------
# 2) Date placeholder on the Handout Master
# ---------------------------------------------------------------------------
$p.HandoutMaster.Shapes.Item(2).TextFrame.TextRange.Text = "9/8/15"

# ---------------------------------------------------------------------------
# 3) Date placeholder on the Notes Master
# ---------------------------------------------------------------------------
$p.NotesMaster.Shapes.Item(2).TextFrame.TextRange.Text = "9/8/15"

# ---------------------------------------------------------------------------
# 4) Slide 1 title: merge "Data " + "Pipeline Demo" runs into one run.
#    (Re-assigning the identical concatenated text is a no-op for the engine,
#    so nudge it through a different string first to force the runs to merge.)
# ---------------------------------------------------------------------------
$titleShape = $p.Slides.Item(1).Shapes.Item(1)
$titleShape.TextFrame.TextRange.Text = "Data Pipeline Demo__tmp__"
$titleShape.TextFrame.TextRange.Text = "Data Pipeline Demo"

# ---------------------------------------------------------------------------
# 5) Slide 14: drop the empty "Text Placeholder 2" (idx=10) shape and reword
#    the subtitle.
# ---------------------------------------------------------------------------
$slide14 = $p.Slides.Item(14)
$slide14.Shapes.Item(2).Delete()
$slide14.Shapes.Item(2).TextFrame.TextRange.Text = "We do Hadoop "

# ---------------------------------------------------------------------------
# 6) Slide 15: drop the empty "Text Placeholder 2" (idx=10) shape.
# ---------------------------------------------------------------------------
$slide15 = $p.Slides.Item(15)
$slide15.Shapes.Item(2).Delete()
